$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data set (rows 2-19), columns A:E ($null = cell left blank)
$data = @(
    @(39400, 2007, -0.7015558851707571, 2008, $null),
    @(39765, 2008, 0.3590181115727509, 2009, $null),
    @(40130, 2009, -0.01587181126743165, 2010, $null),
    @(40494, 2010, -0.02256889165885845, 2011, 0.1850158025574977),
    @(40862, 2011, 0.09611428386597787, 2012, -0.0461580488825053),
    @(41228, 2012, -0.1827723404408288, 2013, -0.1118170657869455),
    @(41592, 2013, -0.001350220946483294, 2014, -0.02003709364999384),
    @(41957, 2014, -0.075754880139145, 2015, -0.2383077634182995),
    @(42321, 2015, -0.5761528471665445, 2016, -0.2318455351884685),
    @(42689, 2016, -0.2011999787958185, 2017, -0.2555440101933648),
    @(43053, 2017, 0.1213692818849532, 2018, -0.06802695342326137),
    @(43418, 2018, 0.1493219406571766, 2019, 0.2006993856294326),
    @(43783, 2019, -0.4278219446121612, 2020, -0.5866687821558636),
    @(44159, 2020, -1.026566979837418, 2021, -2.497698913894009),
    @(44525, 2021, 0.3179894933462268, 2022, 0.1527842396192636),
    @(44890, 2022, 0.463604920919658, 2023, -0.2029074457040814),
    @(45254, 2023, 0.621639092134818, 2024, -0.1003146131857147),
    @(45618, 2024, -0.6768900623516982, 2025, 0.6778422458094902)
)

# Row 19 is brand new; give A19 the same date style as the rest of column A.
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    if ($null -eq $rec[2]) {
        $ws.Cells.Item($row, 3).Value = $null
    } else {
        $ws.Cells.Item($row, 3).Value = $rec[2]
    }
    $ws.Cells.Item($row, 4).Value = $rec[3]
    if ($null -eq $rec[4]) {
        $ws.Cells.Item($row, 5).Value = $null
    } else {
        $ws.Cells.Item($row, 5).Value = $rec[4]
    }
    $row = $row + 1
}
